# Updating parameters used for SOM
# - Change the rows/cols (E/F) grid search values on the "full" sheet
# - Extend rows 14:25 (previously hard-coded/static) to use the same
#   live formulas as rows 2:13 (C = E*F, M = F/E), matching their formatting
# - Refresh the "1.5 ratio .../205 Vesanto nodes" note cells on "full" to the
#   new "1.3 ratio .../175 Vesanto nodes" text
# - Make "full" the active/selected sheet (was "shortened")

$wb = $excel.ActiveWorkbook
$wsFull = $wb.Worksheets.Item("full")
$wsShort = $wb.Worksheets.Item("shortened")

# New rows(E)/cols(F) pairs for rows 2-25 (same 6-row pattern repeated 4x,
# once per nClusters group: 2, 3, 4, 5)
$efPairs = @(
    @(11,15), @(12,15), @(13,15), @(11,16), @(12,16), @(11,17),
    @(11,15), @(12,15), @(13,15), @(11,16), @(12,16), @(11,17),
    @(11,15), @(12,15), @(13,15), @(11,16), @(12,16), @(11,17),
    @(11,15), @(12,15), @(13,15), @(11,16), @(12,16), @(11,17)
)

for ($i = 0; $i -lt $efPairs.Count; $i++) {
    $row = $i + 2
    $wsFull.Cells.Item($row, 5).Value = $efPairs[$i][0]
    $wsFull.Cells.Item($row, 6).Value = $efPairs[$i][1]
}

# Rows 14-25 used to hold static, hand-typed values (with a distinct black-
# font style) instead of the live formulas used by rows 2-13. Bring them in
# line: same "=E*F" / "=F/E" formulas, same (default) formatting.
for ($row = 14; $row -le 25; $row++) {
    $wsFull.Cells.Item($row, 3).Formula = "=E" + $row + "*F" + $row
    $wsFull.Cells.Item($row, 13).Formula = "=F" + $row + "/E" + $row
}

$wsFull.Range("B14:L25").ClearFormats() | Out-Null
$wsFull.Range("M14:M25").ClearFormats() | Out-Null
$wsFull.Range("M14:M25").NumberFormat = "0.0"

# Update the descriptive note cells on "full" (shortened keeps the old text)
$wsFull.Range("O1").Value = "1.3 ratio of cols to rows"
$wsFull.Range("O2").Value = "175 Vesanto nodes"

# "full" becomes the selected/active sheet (was "shortened"); reset its
# selection back to just A1 (was A1:XFD3)
$wsFull.Activate() | Out-Null
$wsFull.Range("A1").Select() | Out-Null
